$d = $word.ActiveDocument

function FindReplaceInRange($range, $findText, $replaceText) {
    $range.Find.Execute($findText, $true, $true, $false, $false, $false, $true, 1, $false, $replaceText, 2) | Out-Null
}

# ---------------------------------------------------------------
# Paragraph 9 (1-indexed): Introduction, first citation paragraph
#   [1] -> [1-4]; sentence split after "...true state of nature";
#   new [2] citation inserted; "Recently, ..." becomes its own sentence;
#   [2-4] -> [5-7]; [5] -> [8]; [6] -> [9]; trailing editorial note added.
# Process back-to-front (end of paragraph first) so that already-settled
# text toward the end of the paragraph is not disturbed by later (leftward) edits.
# ---------------------------------------------------------------

$p9 = $d.Paragraphs(9).Range
FindReplaceInRange $p9 "[6]." "[9]. [not necessarily the best examples]"

$p9 = $d.Paragraphs(9).Range
FindReplaceInRange $p9 "[5], and uncertainty in stock productivity" "[8], and uncertainty in stock productivity"

$p9 = $d.Paragraphs(9).Range
FindReplaceInRange $p9 "[2–4], uncertainty in steepness of the stock-recruit relationship" "[5–7], uncertainty in steepness of the stock-recruit relationship"

$p9 = $d.Paragraphs(9).Range
FindReplaceInRange $p9 "[1]. With simulation, we can evaluate the precision and bias of complex assessment methods in a controlled environment where we know the true state of nature (REFs). Recently, simulation studies have been key to improving strategies for dealing with, for example, time-varying natural mortality" "[1–4]. With simulation, we can evaluate the precision and bias of complex assessment methods in a controlled environment where we know the true state of nature [2]. Recently, simulation studies have been key to improving strategies for dealing with, for example, time-varying natural mortality"

Write-Output "Paragraph 9 done"
Write-Output $d.Paragraphs(9).Range.Text

# ---------------------------------------------------------------
# Paragraph 10 (1-indexed): "Stock Synthesis [7], is a widely-used..."
#   All four "[7]"/"[7,8]" citations renumber to "[10]"/"[10,11]";
#   "models" -> "modeling".
# ---------------------------------------------------------------

$p10 = $d.Paragraphs(10).Range
FindReplaceInRange $p10 "stock assessments, respectively, as of 2012 [7]." "stock assessments, respectively, as of 2012 [10]."

$p10 = $d.Paragraphs(10).Range
FindReplaceInRange $p10 "instead of the model code [7]." "instead of the model code [10]."

$p10 = $d.Paragraphs(10).Range
FindReplaceInRange $p10 "minimally-processed data [7,8]." "minimally-processed data [10,11]."

$p10 = $d.Paragraphs(10).Range
FindReplaceInRange $p10 "population dynamics models using" "population dynamics modeling using"

$p10 = $d.Paragraphs(10).Range
FindReplaceInRange $p10 "Stock Synthesis [7], is a widely-used" "Stock Synthesis [10], is a widely-used"

Write-Output "Paragraph 10 done"
Write-Output $d.Paragraphs(10).Range.Text

# ---------------------------------------------------------------
# Paragraph 11 (1-indexed): "Although SS is increasingly a standard..."
#   [9] -> [12]; "general philosophy" -> "general structure".
# ---------------------------------------------------------------

$p11 = $d.Paragraphs(11).Range
FindReplaceInRange $p11 "outlining the general philosophy of ss3sim" "outlining the general structure of ss3sim"

$p11 = $d.Paragraphs(11).Range
FindReplaceInRange $p11 "[9]" "[12]"

Write-Output "Paragraph 11 done"
Write-Output $d.Paragraphs(11).Range.Text

# ---------------------------------------------------------------
# Paragraph 15 (1-indexed): "Throughout this paper we refer to..."
#   "(OM) to refer" -> "(OM) [13] to refer"
# ---------------------------------------------------------------

$p15 = $d.Paragraphs(15).Range
FindReplaceInRange $p15 "(OM) to refer" "(OM) [13] to refer"

Write-Output "Paragraph 15 done"
Write-Output $d.Paragraphs(15).Range.Text

# ---------------------------------------------------------------
# Paragraph 20 (1-indexed): "Rapid: First, ss3sim relies on SS3..."
#   [10] -> [14]
# ---------------------------------------------------------------

$p20 = $d.Paragraphs(20).Range
FindReplaceInRange $p20 "most rapid and robust optimization software available [10]." "most rapid and robust optimization software available [14]."

Write-Output "Paragraph 20 done"
Write-Output $d.Paragraphs(20).Range.Text
